# daily auto push: 2026-01-07 18:49 UTC
# Insert two new daily records (2026/01/07 23:00 and 2026/01/08 02:00) into the
# log table on Sheet1, just above the existing "2026/12/29" block (old row 602).
# This pushes all subsequent rows down by two, growing the table from
# A1:D643 to A1:D645.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the position right before the old row 602
# (the former first row of the "2026/12/29" block), shifting everything
# below it down by two rows.
$ws.Range("A602:D603").EntireRow.Insert()

# Column A holds dates stored as plain text (e.g. "2026/01/07"), not real
# Excel dates. Force the two new cells to Text format before writing the
# value so Excel does not auto-convert the string into a date serial
# number; then clear the formatting again so the cells end up unstyled,
# matching the rest of the sheet.
$ws.Range("A602:A603").NumberFormat = "@"

$ws.Range("A602").Value = "2026/01/07"
$ws.Range("B602").Value = "水"
$ws.Range("C602").Value = 23
$ws.Range("D602").Value = 201

$ws.Range("A603").Value = "2026/01/08"
$ws.Range("B603").Value = "木"
$ws.Range("C603").Value = 2
$ws.Range("D603").Value = 201

$ws.Range("A602:A603").ClearFormats()
